# "Tried to implement Penality Reward System (unfinished)"
#
# 1) Weekly Sales: append a new weekly row (row 31) one week after the
#    previous last row (45655.99999999999 -> 45662.99999999999), qty 0.
# 2) Daily PO: the single PO data row (row 2) is removed, leaving only
#    the header row.
# 3) Merged (Optional): re-derived from Weekly Sales (cols A/B) and Daily
#    PO (col C, now all zero since the only PO row is gone) - rows 10..31
#    shift to match the new Weekly Sales rows 3..31 (i.e. old row r+1),
#    and the Daily_PO_Qty column collapses to 0 everywhere.
# 4) PO Volume Insights: PO aggregate stats (Total/Average/Max/Min) all
#    drop to 0 now that there is no PO data left.

$wb = $excel.ActiveWorkbook

# --- 1) Weekly Sales: append row 31 ---------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Sales")
$wsWeekly.Range("A31").Value = 45662.99999999999
$wsWeekly.Range("B31").Value = 0
$wsWeekly.Range("A31").NumberFormat = $wsWeekly.Range("A30").NumberFormat

# --- 2) Daily PO: delete the lone data row (row 2) ------------------------
$wsDaily = $wb.Worksheets.Item("Daily PO")
$wsDaily.Rows.Item(2).Delete()

# --- 3) Merged (Optional): rewrite rows 10..31 ----------------------------
$wsMerged = $wb.Worksheets.Item("Merged (Optional)")

$mergedDates = @(45193.99999999999, 45256.99999999999, 45263.99999999999, `
    45270.99999999999, 45277.99999999999, 45291.99999999999, 45298.99999999999, `
    45305.99999999999, 45326.99999999999, 45333.99999999999, 45361.99999999999, `
    45375.99999999999, 45382.99999999999, 45389.99999999999, 45508.99999999999, `
    45536.99999999999, 45606.99999999999, 45634.99999999999, 45641.99999999999, `
    45648.99999999999, 45655.99999999999, 45662.99999999999)

$mergedY = @(0, 0, 0, 4, 1, 1, 2, 1, 0, 1, 1, 1, 1, 1, 0, 2, 1, 2, 1, 1, 0, 0)

for ($i = 0; $i -lt $mergedDates.Length; $i++) {
    $row = 10 + $i
    $wsMerged.Cells.Item($row, 1).Value = $mergedDates[$i]
    $wsMerged.Cells.Item($row, 2).Value = $mergedY[$i]
    $wsMerged.Cells.Item($row, 3).Value = 0
}

# --- 4) PO Volume Insights: aggregates collapse to 0 ----------------------
$wsInsights = $wb.Worksheets.Item("PO Volume Insights")
$wsInsights.Range("A2").Value = 0
$wsInsights.Range("B2").Value = 0
$wsInsights.Range("C2").Value = 0
$wsInsights.Range("D2").Value = 0
